$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was reported for this product/variety. It
# belongs chronologically before the existing row 117 (44939 vs 44827),
# so insert a new row at 117 and shift the rest of the table down.
$ws.Rows.Item(117).Insert()

$ws.Cells.Item(117, 1).Value = 1
$ws.Cells.Item(117, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(117, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(117, 4).Value = 44939
$ws.Cells.Item(117, 5).Value = 15
$ws.Cells.Item(117, 6).Value = 100112021
$ws.Cells.Item(117, 7).Value = "Ají"
$ws.Cells.Item(117, 8).Value = "Cristal"
$ws.Cells.Item(117, 9).Value = "Primera"
$ws.Cells.Item(117, 10).Value = 160
$ws.Cells.Item(117, 11).Value = 9000
$ws.Cells.Item(117, 12).Value = 10000
$ws.Cells.Item(117, 13).Value = 9500
$ws.Cells.Item(117, 14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(117, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(117, 16).Value = 633
$ws.Cells.Item(117, 17).Value = 15
$ws.Cells.Item(117, 18).Value = "Hortaliza"

Write-Output "inserted row 117"
